$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "15-7=",
    "83-19=",
    "65-27=",
    "46+26=",
    "15+68=",
    "90-45=",
    "82-39=",
    "28+3=",
    "60-9=",
    "29+65=",
    "37+9=",
    "18+56=",
    "53-37=",
    "19+14=",
    "87+7=",
    "8+77=",
    "35+46=",
    "21-8=",
    "62+9=",
    "33+9=",
    "39+19=",
    "84-8=",
    "48+5=",
    "45-9=",
    "9+32=",
    "46-27=",
    "44+27=",
    "40-39=",
    "8+18=",
    "4+28=",
    "84+8=",
    "90-21=",
    "7+24=",
    "20-5=",
    "77+5=",
    "55-39=",
    "60-46=",
    "30-16=",
    "6+56=",
    "58-39=",
    "54-36=",
    "78+18=",
    "40-5=",
    "38+47=",
    "70-15=",
    "54+18=",
    "7+24=",
    "49+39=",
    "38-9=",
    "84-29=",
    "37+9=",
    "48+8=",
    "62-57=",
    "20-16=",
    "71-54=",
    "33-17=",
    "14+29=",
    "4+57=",
    "9+28=",
    "49+17=",
    "48+27=",
    "47+17=",
    "64-19=",
    "16+36=",
    "46-18=",
    "16-8=",
    "32+39=",
    "62-19=",
    "19+56=",
    "84-36=",
    "9+56=",
    "28+19=",
    "82-17=",
    "95-19=",
    "73-34=",
    "73-29=",
    "27+59=",
    "75-19=",
    "72-67=",
    "27+56=",
    "51-47=",
    "23+48=",
    "28+18=",
    "30-24=",
    "31-23=",
    "12+9=",
    "3+38=",
    "74-49=",
    "9+67=",
    "9+5=",
    "49+27=",
    "53-34=",
    "14+27=",
    "87-58=",
    "84-6=",
    "47+44=",
    "77+8=",
    "18+77=",
    "19+19=",
    "44+8="
)
$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}
Write-Output "done: $idx cells updated"
